# Update "paises" (countries) COVID-19 stats sheet: new data pull at 19:57
# (previously 18:40), which also changes the row order for a handful of
# countries whose totals crossed each other once re-sorted by "Casos totales".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 19:57"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2699415
$ws.Range("C4").Value = 17604
$ws.Range("D4").Value = 1124142
$ws.Range("E4").Value = 1446196
$ws.Range("G4").Value = 294
$ws.Range("H4").Value = 129077

# Row 7: India
$ws.Range("B7").Value = 585210
$ws.Range("C7").Value = 17674
$ws.Range("D7").Value = 347836
$ws.Range("E7").Value = 219964
$ws.Range("G7").Value = 506
$ws.Range("H7").Value = 17410

# Row 17: Alemania
$ws.Range("B17").Value = 195565
$ws.Range("C17").Value = 173
$ws.Range("E17").Value = 7424

# Row 22: Canada
$ws.Range("B22").Value = 104144
$ws.Range("C22").Value = 226
$ws.Range("D22").Value = 67522
$ws.Range("E22").Value = 28031
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 8591

# Row 51: Irlanda
$ws.Range("B51").Value = 25473
$ws.Range("C51").Value = 11
$ws.Range("E51").Value = 373
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 1736

# Row 53: Israel
$ws.Range("B53").Value = 25041
$ws.Range("C53").Value = 600
$ws.Range("D53").Value = 17318
$ws.Range("E53").Value = 7403

# Row 68: Marruecos
$ws.Range("B68").Value = 12533
$ws.Range("C68").Value = 243
$ws.Range("D68").Value = 8920
$ws.Range("E68").Value = 3385
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 228

# Row 105: Cuba -> Maldivas
$ws.Range("A105").Value = "Maldivas"
$ws.Range("B105").Value = 2361
$ws.Range("C105").Value = 24
$ws.Range("D105").Value = 1944
$ws.Range("E105").Value = 408
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 9

# Row 106: Maldivas -> Cuba
$ws.Range("A106").Value = "Cuba"
$ws.Range("B106").Value = 2341
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 2214
$ws.Range("E106").Value = 41
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 86

# Row 129: Jordania -> Yemen
$ws.Range("A129").Value = "Yemen"
$ws.Range("B129").Value = 1158
$ws.Range("C129").Value = 30
$ws.Range("D129").Value = 488
$ws.Range("E129").Value = 358
$ws.Range("G129").Value = 8
$ws.Range("H129").Value = 312

# Row 130: Yemen -> Jordania
$ws.Range("A130").Value = "Jordania"
$ws.Range("B130").Value = 1132
$ws.Range("C130").Value = 4
$ws.Range("D130").Value = 882
$ws.Range("E130").Value = 241
$ws.Range("H130").Value = 9

# Row 135: Republica de Chipre
$ws.Range("B135").Value = 998
$ws.Range("C135").Value = 2
$ws.Range("E135").Value = 146

# Row 139: Mozambique
$ws.Range("B139").Value = 889
$ws.Range("C139").Value = 6
$ws.Range("D139").Value = 232
$ws.Range("E139").Value = 651

# Row 143: Libia -> Suazilandia
$ws.Range("A143").Value = "Suazilandia"
$ws.Range("B143").Value = 812
$ws.Range("C143").Value = 17
$ws.Range("D143").Value = 408
$ws.Range("E143").Value = 393
$ws.Range("H143").Value = 11

# Row 144: Suazilandia -> Libia
$ws.Range("A144").Value = "Libia"
$ws.Range("B144").Value = 802
$ws.Range("D144").Value = 206
$ws.Range("E144").Value = 573
$ws.Range("H144").Value = 23

# Row 203: Santa Lucia -> Laos
$ws.Range("A203").Value = "Laos"

# Row 204: Laos -> Santa Lucia
$ws.Range("A204").Value = "Santa Lucia"

# Row 209: Groenlandia -> Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"

# Row 210: Islas Malvinas -> Groenlandia
$ws.Range("A210").Value = "Groenlandia"
